$d = $word.ActiveDocument

# Locate the paragraph that currently holds only the page break run
# (it immediately follows the paragraph whose hyperlink text is
# "https://www.youtube.com/watch?v=37UmUAdIHss").
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "37UmUAdIHss") {
        $anchorIndex = $i
    }
}

$pageBreakPara = $d.Paragraphs.Item($anchorIndex + 1)
$insertPos = $pageBreakPara.Range.Start

# Insert all of the new plain text in one shot, right before the
# existing page-break run: "Progress bar:" + manual line break +
# the URL text + a trailing space.
$urlText = "https://www.youtube.com/watch?v=TBGAc1Gj-tM"
$fullText = "Progress bar:" + [char]11 + $urlText + " "

$insertRange = $d.Range($insertPos, $insertPos)
$insertRange.InsertBefore($fullText)

# Now turn just the URL portion of that text into a real hyperlink.
$urlStart = $insertPos + ("Progress bar:" + [char]11).Length
$urlEnd = $urlStart + $urlText.Length
$urlRange = $d.Range($urlStart, $urlEnd)
$d.Hyperlinks.Add($urlRange, $urlText, [Type]::Missing, [Type]::Missing, $urlText) | Out-Null
